$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = "paris"
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = "'"
$ws.Cells.Item(7, 3).ClearFormats()
$ws.Cells.Item(7, 4).Value = "'"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "'"
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(7, 6).Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Cells.Item(7, 7).Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Cells.Item(7, 8).Value = "As you suggested, I did run comparison tests and I will present the results here."
